# Weekly update for Fruta/Hortalizas - Agrícola del Norte S.A. de Arica - Nectarín
# Inserts two new price records (rows 77-78) ahead of the existing data,
# pushing the prior rows 77-106 down to rows 79-108.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 77 (shifts existing rows 77:106 -> 79:108)
$ws.Rows.Item(77).Insert()
$ws.Rows.Item(77).Insert()

# New row 77: Artic Snow / Primera
$ws.Cells.Item(77, 1).Value = 1
$ws.Cells.Item(77, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(77, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(77, 4).Value = 44992
$ws.Cells.Item(77, 5).Value = 15
$ws.Cells.Item(77, 6).Value = "Fruta"
$ws.Cells.Item(77, 7).Value = 100103
$ws.Cells.Item(77, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(77, 9).Value = 100103006
$ws.Cells.Item(77, 10).Value = "Nectarín"
$ws.Cells.Item(77, 11).Value = "Artic Snow"
$ws.Cells.Item(77, 12).Value = "Primera"
$ws.Cells.Item(77, 13).Value = 250
$ws.Cells.Item(77, 14).Value = 20000
$ws.Cells.Item(77, 15).Value = 22000
$ws.Cells.Item(77, 16).Value = 21200
$ws.Cells.Item(77, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(77, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(77, 19).Value = 1178
$ws.Cells.Item(77, 20).Value = 18

# New row 78: August Red / Primera
$ws.Cells.Item(78, 1).Value = 1
$ws.Cells.Item(78, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(78, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(78, 4).Value = 44992
$ws.Cells.Item(78, 5).Value = 15
$ws.Cells.Item(78, 6).Value = "Fruta"
$ws.Cells.Item(78, 7).Value = 100103
$ws.Cells.Item(78, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(78, 9).Value = 100103006
$ws.Cells.Item(78, 10).Value = "Nectarín"
$ws.Cells.Item(78, 11).Value = "August Red"
$ws.Cells.Item(78, 12).Value = "Primera"
$ws.Cells.Item(78, 13).Value = 250
$ws.Cells.Item(78, 14).Value = 20000
$ws.Cells.Item(78, 15).Value = 22000
$ws.Cells.Item(78, 16).Value = 20800
$ws.Cells.Item(78, 17).Value = "`$/bandeja 18 kilos granel"
$ws.Cells.Item(78, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(78, 19).Value = 1156
$ws.Cells.Item(78, 20).Value = 18
